$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H105").Value = 29784
$ws.Range("J105").Value = 29784
$ws.Range("L105").Value = 29784
$ws.Range("N105").Value = -36772
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H137").Value = 1665.4884
$ws.Range("I137").Value = 1164.1562
$ws.Range("J137").Value = 3123.9092
$ws.Range("K137").Value = 3492.4686
$ws.Range("L137").Value = 9371.7276
$ws.Range("M137").Value = -942.4685999999997
$ws.Range("N137").Value = -14471.7276
$ws.Range("H138").Value = 2043160.1
$ws.Range("J138").Value = 3281868.8
$ws.Range("L138").Value = 9845606.399999999
$ws.Range("N138").Value = -9855886.399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 29666.5
$ws.Range("J7").Value = 29666.5
$ws.Range("L7").Value = 29666.5
$ws.Range("N7").Value = -29894.5
$ws.Range("H11").Value = 1675966.6
$ws.Range("J11").Value = 13950
$ws.Range("L11").Value = 13950
$ws.Range("N11").Value = -14238
$ws.Range("H13").Value = 8001600
$ws.Range("I13").Value = 10667400
$ws.Range("J13").Value = 4200
$ws.Range("K13").Value = 10667400
$ws.Range("L13").Value = 4200
$ws.Range("M13").Value = -10667256
$ws.Range("N13").Value = -4488
$ws.Range("H19").Value = 1500
$ws.Range("I19").Value = 1500
$ws.Range("K19").Value = 1500
$ws.Range("M19").Value = -1271
$ws.Range("H61").Value = 2008.5
$ws.Range("I61").Value = 1798.7826
$ws.Range("K61").Value = 1798.7826
$ws.Range("M61").Value = -1586.7826
$ws.Range("H121").Value = 29499.8
$ws.Range("J121").Value = 29499.8
$ws.Range("L121").Value = 29499.8
$ws.Range("N121").Value = -32993.8
$ws.Range("H130").Value = 28714.285
$ws.Range("J130").Value = 28714.285
$ws.Range("L130").Value = 28714.285
$ws.Range("N130").Value = -38754.285
$ws.Range("H136").Value = 2008.5
$ws.Range("I136").Value = 1798.7826
$ws.Range("K136").Value = 5396.3478
$ws.Range("M136").Value = -2846.3478

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 1000000
$ws.Range("J130").Value = 1000000
$ws.Range("L130").Value = 1000000
$ws.Range("N130").Value = -1010040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 3139.5715
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 3139.5715
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 3139.5715
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -3365.5715
$ws.Range("H20").Value = 49999
$ws.Range("J20").Value = 49999
$ws.Range("L20").Value = 49999
$ws.Range("N20").Value = -50471
$ws.Range("H23").Value = 25000
$ws.Range("J23").Value = 25000
$ws.Range("L23").Value = 25000
$ws.Range("N23").Value = -25480
$ws.Range("H27").Value = 25000
$ws.Range("J27").Value = 25000
$ws.Range("L27").Value = 25000
$ws.Range("N27").Value = -25384
$ws.Range("H30").Value = 49999
$ws.Range("J30").Value = 49999
$ws.Range("L30").Value = 49999
$ws.Range("N30").Value = -50181
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H128").Value = 49999
$ws.Range("J128").Value = 49999
$ws.Range("L128").Value = 49999
$ws.Range("N128").Value = -59959
$ws.Range("H129").Value = 32702.762
$ws.Range("J129").Value = 32702.762
$ws.Range("L129").Value = 32702.762
$ws.Range("N129").Value = -42702.762
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 855.24
$ws.Range("I131").Value = 576
$ws.Range("J131").Value = 896.9655
$ws.Range("K131").Value = 1728
$ws.Range("L131").Value = 2690.8965
$ws.Range("M131").Value = 3312
$ws.Range("N131").Value = -12770.8965
$ws.Range("H134").Value = 6256.8335
$ws.Range("I134").Value = 4191.4116
$ws.Range("J134").Value = 7661.32
$ws.Range("K134").Value = 12574.2348
$ws.Range("L134").Value = 22983.96
$ws.Range("M134").Value = -7504.234800000002
$ws.Range("N134").Value = -33123.96

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 4004.75
$ws.Range("I5").Value = 2504.5
$ws.Range("J5").Value = 5505
$ws.Range("K5").Value = 2504.5
$ws.Range("L5").Value = 5505
$ws.Range("M5").Value = -2391.5
$ws.Range("N5").Value = -5731
$ws.Range("H10").Value = 2278
$ws.Range("I10").Value = 499.33334
$ws.Range("J10").Value = 3345.2
$ws.Range("K10").Value = 499.33334
$ws.Range("L10").Value = 3345.2
$ws.Range("M10").Value = -359.33334
$ws.Range("N10").Value = -3625.2
$ws.Range("H17").Value = 11402.667
$ws.Range("I17").Value = 2008
$ws.Range("J17").Value = 16100
$ws.Range("K17").Value = 2008
$ws.Range("L17").Value = 16100
$ws.Range("M17").Value = -1838
$ws.Range("N17").Value = -16440
$ws.Range("H19").Value = 24034.666
$ws.Range("I19").Value = 100
$ws.Range("J19").Value = 36002
$ws.Range("K19").Value = 100
$ws.Range("L19").Value = 36002
$ws.Range("M19").Value = 70
$ws.Range("N19").Value = -36342
$ws.Range("H136").Value = 8131359
$ws.Range("I136").Value = 11906102
$ws.Range("J136").Value = 1143.0769
$ws.Range("K136").Value = 35718306
$ws.Range("L136").Value = 3429.2307
$ws.Range("M136").Value = -35715756
$ws.Range("N136").Value = -8529.2307

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 18322.5
$ws.Range("I10").Value = 9250
$ws.Range("J10").Value = 27395
$ws.Range("K10").Value = 9250
$ws.Range("L10").Value = 27395
$ws.Range("M10").Value = -9081
$ws.Range("N10").Value = -27733
$ws.Range("H12").Value = 13925
$ws.Range("J12").Value = 13925
$ws.Range("L12").Value = 13925
$ws.Range("N12").Value = -14209
$ws.Range("H122").Value = 22846700
$ws.Range("I122").Value = 35729716
$ws.Range("J122").Value = 301422.5
$ws.Range("K122").Value = 107189148
$ws.Range("L122").Value = 904267.5
$ws.Range("M122").Value = -107186698
